# The scraper re-ran and appended 9 more Washington Sundar innings
# (rows 11-19) to the existing table on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: venue, date, result, ownTeam, oppTeam, batsman, totalRuns, totalBalls, total4s, total6s, sr
$newRows = @(
    @(" Dubai (DSC)", " October 25 2020", "Super Kings won by 8 wickets (with 8 balls remaining)", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "5", "2", "1", "0", "250.00"),
    @(" Abu Dhabi", " October 28 2020", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "Washington Sundar ", "10", "6", "1", "0", "166.66"),
    @(" Abu Dhabi", " November 02 2020", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " September 24 2020", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "30", "27", "2", "1", "111.11"),
    @(" Abu Dhabi", " November 06 2020", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "5", "6", "0", "0", "83.33"),
    @(" Sharjah", " October 31 2020", "Sunrisers won by 5 wickets (with 35 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Washington Sundar ", "21", "18", "2", "0", "116.66"),
    @(" Sharjah", " October 15 2020", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "Washington Sundar ", "13", "14", "1", "0", "92.85"),
    @(" Dubai (DSC)", " October 05 2020", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "Washington Sundar ", "17", "11", "3", "0", "154.54"),
    @(" Dubai (DSC)", " October 10 2020", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "Washington Sundar ", "10", "10", "0", "1", "100.00")
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $value = $rowValues[$c - 1]
        if ($c -ge 7) {
            # totalRuns / totalBalls / total4s / total6s / sr are numeric-looking
            # strings that must stay stored as text (matches the source sheet,
            # which has numberStoredAsText warnings suppressed for these cols).
            $ws.Cells.Item($r, $c).Value = "'" + $value
        } else {
            $ws.Cells.Item($r, $c).Value = $value
        }
    }
}
